# Auto-generated edit script: updates recalculated profit-table values
# across the 8 crafting-class worksheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 5262150.5
$ws.Range("J17").Value = 5525233
$ws.Range("L17").Value = 16575699
$ws.Range("N17").Value = -16576035

$ws.Range("H62").Value = 1301
$ws.Range("I62").Value = 1301.6666
$ws.Range("J62").Value = 1300
$ws.Range("K62").Value = 1301.6666
$ws.Range("L62").Value = 1300
$ws.Range("M62").Value = -677.6666
$ws.Range("N62").Value = -2548

$ws.Range("H65").Value = 1301
$ws.Range("I65").Value = 1301.6666
$ws.Range("J65").Value = 1300
$ws.Range("K65").Value = 6508.333000000001
$ws.Range("L65").Value = 6500
$ws.Range("M65").Value = -3388.333000000001
$ws.Range("N65").Value = -12740

$ws.Range("H127").Value = 1369.8206
$ws.Range("I127").Value = 479
$ws.Range("J127").Value = 1531.7878
$ws.Range("K127").Value = 1437
$ws.Range("L127").Value = 4595.3634
$ws.Range("M127").Value = 3523
$ws.Range("N127").Value = -14515.3634

$ws.Range("H129").Value = 1777.1428
$ws.Range("J129").Value = 1808.8235
$ws.Range("L129").Value = 5426.470499999999
$ws.Range("N129").Value = -15426.4705

$ws.Range("H132").Value = 4466511.5
$ws.Range("I132").Value = 2129.02
$ws.Range("J132").Value = 41669700
$ws.Range("K132").Value = 6387.059999999999
$ws.Range("L132").Value = 125009100
$ws.Range("M132").Value = -3857.059999999999
$ws.Range("N132").Value = -125014160

$ws.Range("H137").Value = 9868.161
$ws.Range("I137").Value = 895.45
$ws.Range("J137").Value = 26182.182
$ws.Range("K137").Value = 2686.35
$ws.Range("L137").Value = 78546.546
$ws.Range("M137").Value = -136.3500000000004
$ws.Range("N137").Value = -83646.546

$ws.Range("H138").Value = 3475452.2
$ws.Range("I138").Value = 8773829
$ws.Range("J138").Value = 4102.224
$ws.Range("K138").Value = 26321487
$ws.Range("L138").Value = 12306.672
$ws.Range("M138").Value = -26316347
$ws.Range("N138").Value = -22586.672

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1849.5
$ws.Range("I2").Value = 1424.5
$ws.Range("J2").Value = 2274.5
$ws.Range("K2").Value = 1424.5
$ws.Range("L2").Value = 2274.5
$ws.Range("M2").Value = -1311.5
$ws.Range("N2").Value = -2500.5

$ws.Range("H32").Value = 16697.73
$ws.Range("I32").Value = 16620.307
$ws.Range("J32").Value = 18633.334
$ws.Range("K32").Value = 16620.307
$ws.Range("L32").Value = 18633.334
$ws.Range("M32").Value = -16333.307
$ws.Range("N32").Value = -19207.334

$ws.Range("H61").Value = 1636.4445
$ws.Range("I61").Value = 1643.4839
$ws.Range("K61").Value = 1643.4839
$ws.Range("M61").Value = -1431.4839

$ws.Range("H102").Value = 2108.7273
$ws.Range("I102").Value = 1778.75
$ws.Range("J102").Value = 2988.6667
$ws.Range("K102").Value = 1778.75
$ws.Range("L102").Value = 2988.6667
$ws.Range("M102").Value = -156.75
$ws.Range("N102").Value = -6232.6667

$ws.Range("H116").Value = 1849.5
$ws.Range("I116").Value = 1424.5
$ws.Range("J116").Value = 2274.5
$ws.Range("K116").Value = 1424.5
$ws.Range("L116").Value = 2274.5
$ws.Range("M116").Value = 869.5
$ws.Range("N116").Value = -6862.5

$ws.Range("H136").Value = 1636.4445
$ws.Range("I136").Value = 1643.4839
$ws.Range("K136").Value = 4930.4517
$ws.Range("M136").Value = -2380.4517

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1849.5
$ws.Range("I3").Value = 1424.5
$ws.Range("J3").Value = 2274.5
$ws.Range("K3").Value = 1424.5
$ws.Range("L3").Value = 2274.5
$ws.Range("M3").Value = -1310.5
$ws.Range("N3").Value = -2502.5

$ws.Range("H59").Value = 87950
$ws.Range("J59").Value = 87950
$ws.Range("L59").Value = 87950
$ws.Range("N59").Value = -89644

$ws.Range("H96").Value = 16305.714
$ws.Range("I96").Value = 8856.75
$ws.Range("J96").Value = 26237.666
$ws.Range("K96").Value = 8856.75
$ws.Range("L96").Value = 26237.666
$ws.Range("M96").Value = -6110.75
$ws.Range("N96").Value = -31729.666

$ws.Range("H99").Value = 895
$ws.Range("I99").Value = 895
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 895
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = 603
$ws.Range("N99").ClearContents()

$ws.Range("H105").Value = 4754.125
$ws.Range("I105").Value = 3733
$ws.Range("K105").Value = 3733
$ws.Range("M105").Value = -1986

$ws.Range("H139").Value = 20000
$ws.Range("J139").Value = 20000
$ws.Range("L139").Value = 20000
$ws.Range("N139").Value = -30280

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5878.25
$ws.Range("I31").Value = 2146.7368
$ws.Range("J31").Value = 7794.4326
$ws.Range("K31").Value = 2146.7368
$ws.Range("L31").Value = 7794.4326
$ws.Range("M31").Value = -1851.7368
$ws.Range("N31").Value = -8384.4326

$ws.Range("H34").Value = 5878.25
$ws.Range("I34").Value = 2146.7368
$ws.Range("J34").Value = 7794.4326
$ws.Range("K34").Value = 2146.7368
$ws.Range("L34").Value = 7794.4326
$ws.Range("M34").Value = -1944.7368
$ws.Range("N34").Value = -8198.4326

$ws.Range("H122").Value = 2411.6667
$ws.Range("I122").Value = 2806.25
$ws.Range("J122").Value = 1622.5
$ws.Range("K122").Value = 8418.75
$ws.Range("L122").Value = 4867.5
$ws.Range("M122").Value = -5968.75
$ws.Range("N122").Value = -9767.5

$ws.Range("H134").Value = 581761.2
$ws.Range("I134").Value = 1088.9744
$ws.Range("J134").Value = 11904869
$ws.Range("K134").Value = 3266.9232
$ws.Range("L134").Value = 35714607
$ws.Range("M134").Value = -731.9232000000002
$ws.Range("N134").Value = -35719677

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 3697.7317
$ws.Range("I137").Value = 3929
$ws.Range("J137").Value = 3623.1292
$ws.Range("K137").Value = 11787
$ws.Range("L137").Value = 10869.3876
$ws.Range("M137").Value = -6687
$ws.Range("N137").Value = -21069.3876

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H123").Value = 24774.39
$ws.Range("J123").Value = 24774.39
$ws.Range("L123").Value = 24774.39
$ws.Range("N123").Value = -29674.39

$ws.Range("H126").Value = 4587.875
$ws.Range("I126").Value = 3742
$ws.Range("J126").Value = 4972.364
$ws.Range("K126").Value = 11226
$ws.Range("L126").Value = 14917.092
$ws.Range("M126").Value = -8756
$ws.Range("N126").Value = -19857.092

$ws.Range("H132").Value = 6604.577
$ws.Range("I132").Value = 7591.2
$ws.Range("J132").Value = 3315.8333
$ws.Range("K132").Value = 22773.6
$ws.Range("L132").Value = 9947.499899999999
$ws.Range("M132").Value = -20243.6
$ws.Range("N132").Value = -15007.4999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2242.6667
$ws.Range("I68").Value = 2039.8
$ws.Range("J68").Value = 2749.8333
$ws.Range("K68").Value = 2039.8
$ws.Range("L68").Value = 2749.8333
$ws.Range("M68").Value = -1290.8
$ws.Range("N68").Value = -4247.8333

$ws.Range("H71").Value = 2242.6667
$ws.Range("I71").Value = 2039.8
$ws.Range("J71").Value = 2749.8333
$ws.Range("K71").Value = 10199
$ws.Range("L71").Value = 13749.1665
$ws.Range("M71").Value = -6455
$ws.Range("N71").Value = -21237.1665

$ws.Range("H136").Value = 4848.476
$ws.Range("I136").Value = 2188.3125
$ws.Range("J136").Value = 13361
$ws.Range("K136").Value = 6564.9375
$ws.Range("L136").Value = 40083
$ws.Range("M136").Value = -4014.9375
$ws.Range("N136").Value = -45183

$ws.Range("H139").Value = 100000
$ws.Range("J139").Value = 100000
$ws.Range("L139").Value = 100000
$ws.Range("N139").Value = -110280

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2168.125
$ws.Range("I122").Value = 2563.7646
$ws.Range("J122").Value = 1207.2858
$ws.Range("K122").Value = 7691.293799999999
$ws.Range("L122").Value = 3621.8574
$ws.Range("M122").Value = -5241.293799999999
$ws.Range("N122").Value = -8521.857400000001

$ws.Range("H126").Value = 1902.2439
$ws.Range("I126").Value = 1737.9615
$ws.Range("J126").Value = 2187
$ws.Range("K126").Value = 5213.8845
$ws.Range("L126").Value = 6561
$ws.Range("M126").Value = -2743.8845
$ws.Range("N126").Value = -11501

$ws.Range("H132").Value = 1203.5186
$ws.Range("I132").Value = 903.72
$ws.Range("K132").Value = 2711.16
$ws.Range("M132").Value = -181.1599999999999

$ws.Range("H136").Value = 1399.6
$ws.Range("I136").Value = 716.2222
$ws.Range("J136").Value = 7550
$ws.Range("K136").Value = 2148.6666
$ws.Range("L136").Value = 22650
$ws.Range("M136").Value = 401.3334
$ws.Range("N136").Value = -27750

$ws.Range("H139").Value = 43957.5
$ws.Range("J139").Value = 43957.5
$ws.Range("L139").Value = 43957.5
$ws.Range("N139").Value = -54237.5
